# feat(arsenal): Registrar CZ P-10 C para Ricardo Soberanis (Credencial 230)
# Inserts a new row at 283 for a new firearm entry belonging to
# RICARDO ANTONIO SOBERANIS GAMBOA (credencial 230): a CZ P-10 C pistol,
# shifting the existing rows 283-287 down to 284-288.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 283 (existing rows 283:287 shift to 284:288)
$ws.Rows.Item(283).Insert()

# --- Populate the new row 283 ---
$ws.Range("A283").Value = 738
$ws.Range("B283").Value = "CALLE 50 No. 531-E x 69 y 71, CENTRO, 97000 MÉRIDA, YUC."
$ws.Range("C283").Value = 230
$ws.Range("D283").Value = "RICARDO ANTONIO SOBERANIS GAMBOA"
$ws.Range("E283").Value = "SOGR701015HYNBMC04"
$ws.Range("F283").Value = "''9993437376"
$ws.Range("G283").Value = "rsoberanis11@hotmail.com"
$ws.Range("H283").Value = 45808
$ws.Range("I283").Value = "Calle 23 S/N T.C. 50641"
$ws.Range("J283").Value = "Colonia San Antonio Hool"
$ws.Range("K283").Value = "Mérida"
$ws.Range("L283").Value = "YUCATÁN"
$ws.Range("M283").Value = "''97302"
$ws.Range("N283").Value = "PISTOLA"
$ws.Range("O283").Value = ".40 S&W"
$ws.Range("P283").Value = "CZ"
$ws.Range("Q283").Value = "P-10 C"
$ws.Range("R283").Value = "EP29710"
$ws.Range("S283").Value = "A3912487"

# The row-insert carries the formatting of the row above (280-282, the other
# Ricardo Soberanis rows) into the new row's F/M/Q/R/H cells. The phone
# number (F) and postal code (M) cells on this new row should end up with no
# explicit cell style (matching the target data), and likewise the new gun
# columns Q/R (model/matricula) should not inherit the "text" style either.
# Clear AFTER the values are set (clearing then re-setting the value would
# regenerate a quote-prefix style instead of leaving the cell unstyled).
$ws.Range("F283").ClearFormats()
$ws.Range("M283").ClearFormats()
$ws.Range("Q283").ClearFormats()
$ws.Range("R283").ClearFormats()

# H283 (FECHA ALTA) keeps the date number format used throughout column H.
$ws.Range("H283").NumberFormat = "yyyy-mm-dd h:mm:ss"

Write-Host "Row 283 inserted and populated for RICARDO ANTONIO SOBERANIS GAMBOA (CZ P-10 C)."
